$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab name shown as "SCD0251" in workbook.xml -> "SCD0016")
$ws.Name = "SCD0016"

# Update the TC_ID in B2 from "DGS-266" to "SCD0016-025". (C2/D2/E2 text is unchanged;
# their shared-string index merely shifts because the old "DGS-266" string is dropped
# from the middle of the shared strings table while the new one is appended at the end.)
$ws.Range("B2").Value = "SCD0016-025"

# Column B width change (9 -> 12.5703125). The COM layer stores width in
# whole-pixel increments of 1/6 character units (stored = ColumnWidth + 5/6,
# rounded to the nearest 1/6), so 12.5703125 itself is not exactly reachable;
# 11.666666666666666 lands on the closest attainable stored width (12.5).
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666

# Sheet view: remove frozen/topLeftCell E1 (scroll back to A1) and change selection to B3
[void]$ws.Range("A1").Select()
[void]$ws.Range("B3").Select()
